$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 668.8
$ws.Range("I9").Value = 804.9167
$ws.Range("K9").Value = 804.9167
$ws.Range("M9").Value = -635.9167
$ws.Range("H17").Value = 2145.1853
$ws.Range("J17").Value = 2177.6924
$ws.Range("L17").Value = 6533.0772
$ws.Range("N17").Value = -6869.0772
$ws.Range("H92").Value = 174.16667
$ws.Range("I92").Value = 119.4
$ws.Range("K92").Value = 119.4
$ws.Range("M92").Value = 1128.6
$ws.Range("H132").Value = 3147.9565
$ws.Range("I132").Value = 3287.4375
$ws.Range("K132").Value = 9862.3125
$ws.Range("M132").Value = -7332.3125
$ws.Range("H137").Value = 775014.9399999999
$ws.Range("I137").Value = 627502.25
$ws.Range("J137").Value = 1011035.3
$ws.Range("K137").Value = 1882506.75
$ws.Range("L137").Value = 3033105.9
$ws.Range("M137").Value = -1879956.75
$ws.Range("N137").Value = -3038205.9
$ws.Range("H138").Value = 7529.2407
$ws.Range("I138").Value = 3800.1
$ws.Range("K138").Value = 11400.3
$ws.Range("M138").Value = -6260.299999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4806.8335
$ws.Range("I32").Value = 4210.339
$ws.Range("K32").Value = 4210.339
$ws.Range("M32").Value = -3923.339
$ws.Range("H43").Value = 24392
$ws.Range("J43").Value = 24638.5
$ws.Range("L43").Value = 24638.5
$ws.Range("N43").Value = -25264.5
$ws.Range("H61").Value = 4737.25
$ws.Range("I61").Value = 2749.5
$ws.Range("K61").Value = 2749.5
$ws.Range("M61").Value = -2537.5
$ws.Range("H74").Value = 3250
$ws.Range("I74").Value = 3250
$ws.Range("K74").Value = 3250
$ws.Range("M74").Value = -2376
$ws.Range("H77").Value = 3250
$ws.Range("I77").Value = 3250
$ws.Range("K77").Value = 16250
$ws.Range("M77").Value = -11882
$ws.Range("H132").Value = 3147.9343
$ws.Range("I132").Value = 2356.204
$ws.Range("K132").Value = 7068.612000000001
$ws.Range("M132").Value = -4538.612000000001
$ws.Range("H136").Value = 4737.25
$ws.Range("I136").Value = 2749.5
$ws.Range("K136").Value = 8248.5
$ws.Range("M136").Value = -5698.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4432.5
$ws.Range("I99").Value = 3774.3333
$ws.Range("K99").Value = 3774.3333
$ws.Range("M99").Value = -2276.3333
$ws.Range("H122").Value = 73997
$ws.Range("J122").Value = 73997
$ws.Range("L122").Value = 73997
$ws.Range("N122").Value = -83797
$ws.Range("H134").Value = 50811.617
$ws.Range("I134").Value = 2749.7144
$ws.Range("J134").Value = 146935.42
$ws.Range("K134").Value = 8249.143199999999
$ws.Range("L134").Value = 440806.26
$ws.Range("M134").Value = -5714.143199999999
$ws.Range("N134").Value = -445876.26

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35066.88
$ws.Range("I31").Value = 1538.0588
$ws.Range("J31").Value = 70691.25
$ws.Range("K31").Value = 1538.0588
$ws.Range("L31").Value = 70691.25
$ws.Range("M31").Value = -1243.0588
$ws.Range("N31").Value = -71281.25
$ws.Range("H34").Value = 35066.88
$ws.Range("I34").Value = 1538.0588
$ws.Range("J34").Value = 70691.25
$ws.Range("K34").Value = 1538.0588
$ws.Range("L34").Value = 70691.25
$ws.Range("M34").Value = -1336.0588
$ws.Range("N34").Value = -71095.25
$ws.Range("H58").Value = 390632.3
$ws.Range("I58").Value = 593214.1
$ws.Range("K58").Value = 593214.1
$ws.Range("M58").Value = -593011.1
$ws.Range("H62").Value = 3502.5
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3502.5
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H132").Value = 4354.2646
$ws.Range("I132").Value = 3999.28
$ws.Range("J132").Value = 5340.3335
$ws.Range("K132").Value = 11997.84
$ws.Range("L132").Value = 16021.0005
$ws.Range("M132").Value = -9467.84
$ws.Range("N132").Value = -21081.0005
$ws.Range("H136").Value = 390632.3
$ws.Range("I136").Value = 593214.1
$ws.Range("K136").Value = 1779642.3
$ws.Range("M136").Value = -1777092.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11566899

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 410832.72
$ws.Range("I132").Value = 628268.1
$ws.Range("K132").Value = 1884804.3
$ws.Range("M132").Value = -1882274.3

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 965
$ws.Range("I22").Value = 947.5
$ws.Range("K22").Value = 947.5
$ws.Range("M22").Value = -652.5
$ws.Range("H27").Value = 965
$ws.Range("I27").Value = 947.5
$ws.Range("K27").Value = 947.5
$ws.Range("M27").Value = -840.5
$ws.Range("H46").Value = 5684.1055
$ws.Range("I46").Value = 6454.364
$ws.Range("K46").Value = 6454.364
$ws.Range("M46").Value = -6266.364
$ws.Range("H61").Value = 6094.3125
$ws.Range("I61").Value = 5761.1
$ws.Range("J61").Value = 6649.6665
$ws.Range("K61").Value = 5761.1
$ws.Range("L61").Value = 6649.6665
$ws.Range("M61").Value = -5559.1
$ws.Range("N61").Value = -7053.6665
$ws.Range("H93").Value = 71431210
$ws.Range("J93").Value = 3855.5
$ws.Range("L93").Value = 3855.5
$ws.Range("N93").Value = -6351.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 6094.3125
$ws.Range("I113").Value = 5761.1
$ws.Range("J113").Value = 6649.6665
$ws.Range("K113").Value = 5761.1
$ws.Range("L113").Value = 6649.6665
$ws.Range("M113").Value = -3591.1
$ws.Range("N113").Value = -10989.6665
$ws.Range("H130").Value = 90000
$ws.Range("J130").Value = 90000
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040
$ws.Range("H132").Value = 5577.273
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5577.273
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 16731.819
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -21791.819
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1741.1389
$ws.Range("J126").Value = 2830.111
$ws.Range("L126").Value = 8490.332999999999
$ws.Range("N126").Value = -13430.333
$ws.Range("H132").Value = 55422.95
$ws.Range("I132").Value = 5244.4165
$ws.Range("K132").Value = 15733.2495
$ws.Range("M132").Value = -13203.2495
